$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 86: a new date label in column A, and a reused price-text string in column B.
# Column A: enter the date-looking text via a formula so Excel's "looks like a date"
# auto-conversion never kicks in, then collapse the formula down to a plain literal
# value (matches how the other date cells in the sheet are stored as shared strings).
$ws.Range("A86").Formula = "=""11-12-2025"""
$ws.Range("A86").Copy()
$ws.Range("A86").PasteSpecial(-4163)  # xlPasteValues

# Column B reuses the same shared text as row 77 (same gold-price paragraph).
$ws.Range("B86").Value2 = $ws.Range("B77").Value2

$excel.CutCopyMode = 0
